$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, avoiding Excel's
# auto-conversion of numeric-looking strings into numbers, while
# keeping the cell's original (default) style/format untouched.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '64.910.89'
$ws.Range('E2').Value = '  +0.17%  '
Set-TextValue $ws.Range('D3') '3.481.90'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '577.64'
$ws.Range('E5').Value = '  +0.33%  '
Set-TextValue $ws.Range('D6') '161.45'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue $ws.Range('D8') '3.480.18'
$ws.Range('E8').Value = '  +0.73%  '
Set-TextValue $ws.Range('D9') '0.579'
$ws.Range('E9').Value = '  -7.55%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('E14').Value = '  -0.19%  '
Set-TextValue $ws.Range('D15') '27.68'
$ws.Range('E15').Value = '  -2.14%  '
Set-TextValue $ws.Range('D16') '0.0000176'
$ws.Range('E16').Value = '  -8.20%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '3.567.69'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D18') '64.970.44'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').Value = '  -3.49%  '
Set-TextValue $ws.Range('D20') '13.89'
$ws.Range('E20').Value = '  -3.25%  '
Set-TextValue $ws.Range('D21') '382.76'
$ws.Range('E21').Value = '  +0.60%  '
Set-TextValue $ws.Range('D22') '8.00'
$ws.Range('E22').Value = '  -1.14%  '
Set-TextValue $ws.Range('D23') '72.80'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  -3.25%  '
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('E27').Value = '  -1.50%  '
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -3.80%  '
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('E32').Value = '  -0.90%  '
Set-TextValue $ws.Range('D33') '23.45'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('E34').Value = '  -2.11%  '
$ws.Range('E35').Value = '  -1.31%  '
Set-TextValue $ws.Range('D36') '161.55'
$ws.Range('E36').Value = '  +0.50%  '
Set-TextValue $ws.Range('D37') '1.89'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D38') '27.08'
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D39') '0.0755'
$ws.Range('E39').Value = '  -2.61%  '
Set-TextValue $ws.Range('D40') '2.879.33'
$ws.Range('E40').Value = '  -2.14%  '
Set-TextValue $ws.Range('D41') '0.815'
$ws.Range('E41').Value = '  +5.51%  '
$ws.Range('E42').Value = '  -1.57%  '
Set-TextValue $ws.Range('D43') '6.57'
$ws.Range('E43').Value = '  -1.95%  '
Set-TextValue $ws.Range('D44') '43.02'
$ws.Range('E44').Value = '  +0.55%  '
Set-TextValue $ws.Range('D45') '26.03'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('E46').Value = '  -2.70%  '
$ws.Range('E47').Value = '  +11.97%  '
Set-TextValue $ws.Range('D48') '331.93'
$ws.Range('E48').Value = '  +3.04%  '
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('E50').Value = '  -3.05%  '
Set-TextValue $ws.Range('D51') '6.49'
$ws.Range('E51').Value = '  -1.83%  '
